$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.88
$ws.Range("I2").Value = 2.75
$ws.Range("J2").Value = 3.75
$ws.Range("M2").Value = 1.11
$ws.Range("N2").Value = 6.5
$ws.Range("O2").Value = 1.5
$ws.Range("P2").Value = 2.5
$ws.Range("U2").Value = 2.1
$ws.Range("V2").Value = 1.67
$ws.Range("W2").Value = 7
$ws.Range("AE2").Value = 17
$ws.Range("AH2").Value = 7
$ws.Range("AN2").Value = 4.75
$ws.Range("BB2").Value = 301
$ws.Range("G3").Value = 2
$ws.Range("H3").Value = 3.1
$ws.Range("I3").Value = 4
$ws.Range("J3").Value = 2.75
$ws.Range("K3").Value = 1.95
$ws.Range("L3").Value = 4.75
$ws.Range("Q3").Value = 2.4
$ws.Range("R3").Value = 1.53
$ws.Range("U3").Value = 2.1
$ws.Range("V3").Value = 1.67
$ws.Range("W3").Value = 6
$ws.Range("X3").Value = 8.5
$ws.Range("Y3").Value = 9.5
$ws.Range("Z3").Value = 17
$ws.Range("AA3").Value = 19
$ws.Range("AD3").Value = 6
$ws.Range("AE3").Value = 19
$ws.Range("AF3").Value = 67
$ws.Range("AH3").Value = 9
$ws.Range("AI3").Value = 19
$ws.Range("AJ3").Value = 15
$ws.Range("AK3").Value = 41
$ws.Range("AN3").Value = 3.75
$ws.Range("AO3").Value = 12
$ws.Range("AQ3").Value = 41
$ws.Range("AU3").Value = 9
$ws.Range("AV3").Value = 67
$ws.Range("AW3").Value = 5.5
$ws.Range("AX3").Value = 23
$ws.Range("AY3").Value = 34
$ws.Range("AZ3").Value = 81
$ws.Range("BA3").Value = 126
$ws.Range("BB3").Value = 351
$ws.Range("G4").Value = 5.5
$ws.Range("H4").Value = 3.4
$ws.Range("I4").Value = 1.7
$ws.Range("J4").Value = 6
$ws.Range("L4").Value = 2.4
$ws.Range("Q4").Value = 2.35
$ws.Range("R4").Value = 1.57
$ws.Range("W4").Value = 11
$ws.Range("X4").Value = 26
$ws.Range("Y4").Value = 19
$ws.Range("AQ4").Value = 126
$ws.Range("AW4").Value = 3.5
$ws.Range("AX4").Value = 9.5
$ws.Range("AZ4").Value = 34
$ws.Range("H5").Value = 4.33
$ws.Range("I5").Value = 8.5
$ws.Range("M5").Value = 1.06
$ws.Range("N5").Value = 10
$ws.Range("W5").Value = 5.5
$ws.Range("Y5").Value = 9.5
$ws.Range("Z5").Value = 8
$ws.Range("AC5").Value = 9
$ws.Range("AD5").Value = 8.5
$ws.Range("AQ5").Value = 21
$ws.Range("AU5").Value = 11
$ws.Range("AW5").Value = 9.5
$ws.Range("AZ5").Value = 251
$ws.Range("G6").Value = 1.9
$ws.Range("H6").Value = 3.2
$ws.Range("I6").Value = 4.33
$ws.Range("J6").Value = 2.63
$ws.Range("L6").Value = 5
$ws.Range("Q6").Value = 2.35
$ws.Range("R6").Value = 1.57
$ws.Range("X6").Value = 8
$ws.Range("Y6").Value = 9
$ws.Range("Z6").Value = 15
$ws.Range("AD6").Value = 6.5
$ws.Range("AE6").Value = 19
$ws.Range("AH6").Value = 10
$ws.Range("AI6").Value = 21
$ws.Range("AK6").Value = 51
$ws.Range("AO6").Value = 11
$ws.Range("AS6").Value = 201
$ws.Range("AW6").Value = 6
$ws.Range("AX6").Value = 26
$ws.Range("AY6").Value = 41
$ws.Range("AZ6").Value = 101
$ws.Range("G7").Value = 1.5
$ws.Range("H7").Value = 3.8
$ws.Range("I7").Value = 6.5
$ws.Range("J7").Value = 2.1
$ws.Range("K7").Value = 2.2
$ws.Range("L7").Value = 7
$ws.Range("M7").Value = 1.06
$ws.Range("N7").Value = 10
$ws.Range("Q7").Value = 2.1
$ws.Range("R7").Value = 1.7
$ws.Range("S7").Value = 1.44
$ws.Range("T7").Value = 2.63
$ws.Range("U7").Value = 2.2
$ws.Range("V7").Value = 1.62
$ws.Range("X7").Value = 6.5
$ws.Range("Z7").Value = 10
$ws.Range("AC7").Value = 8.5
$ws.Range("AD7").Value = 7.5
$ws.Range("AF7").Value = 81
$ws.Range("AH7").Value = 13
$ws.Range("AI7").Value = 34
$ws.Range("AJ7").Value = 21
$ws.Range("AK7").Value = 81
$ws.Range("AL7").Value = 51
$ws.Range("AN7").Value = 3.25
$ws.Range("AO7").Value = 8
$ws.Range("AQ7").Value = 26
$ws.Range("AT7").Value = 2.63
$ws.Range("AU7").Value = 10
$ws.Range("AV7").Value = 81
$ws.Range("AW7").Value = 8
$ws.Range("AX7").Value = 41
$ws.Range("AZ7").Value = 151
$ws.Range("BA7").Value = 201
$ws.Range("G8").Value = 1.4
$ws.Range("H8").Value = 4.2
$ws.Range("I8").Value = 8.5
$ws.Range("J8").Value = 1.95
$ws.Range("L8").Value = 8
$ws.Range("Q8").Value = 2.03
$ws.Range("R8").Value = 1.83
$ws.Range("W8").Value = 5.5
$ws.Range("AD8").Value = 8.5
$ws.Range("AH8").Value = 17
$ws.Range("AI8").Value = 41
$ws.Range("AK8").Value = 101
$ws.Range("AL8").Value = 67
$ws.Range("AM8").Value = 67
$ws.Range("BA8").Value = 251
$ws.Range("G9").Value = 1.38
$ws.Range("H9").Value = 4.33
$ws.Range("I9").Value = 8
$ws.Range("K9").Value = 2.3
$ws.Range("L9").Value = 8
$ws.Range("O9").Value = 1.3
$ws.Range("P9").Value = 3.4
$ws.Range("Q9").Value = 2
$ws.Range("R9").Value = 1.85
$ws.Range("Z9").Value = 8.5
$ws.Range("AK9").Value = 101
$ws.Range("AL9").Value = 67
$ws.Range("AW9").Value = 9
$ws.Range("G10").Value = 1.85
$ws.Range("H10").Value = 2.9
$ws.Range("I10").Value = 5.25
$ws.Range("J10").Value = 2.75
$ws.Range("L10").Value = 6
$ws.Range("M10").Value = 1.13
$ws.Range("N10").Value = 6
$ws.Range("O10").Value = 1.57
$ws.Range("P10").Value = 2.25
$ws.Range("Q10").Value = 2.88
$ws.Range("R10").Value = 1.4
$ws.Range("X10").Value = 7
$ws.Range("Y10").Value = 10
$ws.Range("Z10").Value = 15
$ws.Range("AA10").Value = 21
$ws.Range("AC10").Value = 5.5
$ws.Range("AD10").Value = 6.5
$ws.Range("AH10").Value = 9.5
$ws.Range("AI10").Value = 23
$ws.Range("AN10").Value = 3.6
$ws.Range("AO10").Value = 11
$ws.Range("AP10").Value = 29
$ws.Range("AR10").Value = 81
$ws.Range("AW10").Value = 6.5
$ws.Range("AX10").Value = 34
$ws.Range("G11").Value = 1.5
$ws.Range("K11").Value = 2.2
$ws.Range("AC11").Value = 8.5
$ws.Range("AD11").Value = 7.5
$ws.Range("AE11").Value = 21
$ws.Range("AM11").Value = 51
$ws.Range("AO11").Value = 8
$ws.Range("AQ11").Value = 26
$ws.Range("N12").Value = 5.95
$ws.Range("H13").Value = 3.1
$ws.Range("I13").Value = 4.05
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 4.4
$ws.Range("M13").Value = 1.01
$ws.Range("N13").Value = 9.449999999999999
$ws.Range("O13").Value = 1.28
$ws.Range("P13").Value = 3.05
$ws.Range("Q13").Value = 1.88
$ws.Range("R13").Value = 1.83
$ws.Range("S13").Value = 1.42
$ws.Range("T13").Value = 2.47
$ws.Range("U13").Value = 1.65
$ws.Range("V13").Value = 1.98
$ws.Range("Z13").Value = 18
$ws.Range("AC13").Value = 9.25
$ws.Range("AD13").Value = 6.1
$ws.Range("AF13").Value = 55
$ws.Range("AG13").Value = 400
$ws.Range("AH13").Value = 12
$ws.Range("AJ13").Value = 13
$ws.Range("AL13").Value = 37
$ws.Range("AM13").Value = 37
$ws.Range("AN13").Value = 3.8
$ws.Range("AS13").Value = 250
$ws.Range("AT13").Value = 2.45
$ws.Range("AU13").Value = 6.7
$ws.Range("AV13").Value = 60
$ws.Range("AX13").Value = 23
$ws.Range("AY13").Value = 27
$ws.Range("BA13").Value = 150
$ws.Range("BB13").Value = 350
$ws.Range("M14").Value = 1.08
$ws.Range("N14").Value = 8
$ws.Range("Q15").Value = 2.08
$ws.Range("R15").Value = 1.73
$ws.Range("I16").Value = 1.57
$ws.Range("M16").Value = 1.05
$ws.Range("N16").Value = 11
$ws.Range("O16").Value = 1.3
$ws.Range("P16").Value = 3.4
$ws.Range("Q16").Value = 2.03
$ws.Range("R16").Value = 1.83
$ws.Range("U16").Value = 2.1
$ws.Range("V16").Value = 1.67
$ws.Range("X16").Value = 29
$ws.Range("AE16").Value = 21
$ws.Range("AW16").Value = 3.4
$ws.Range("AX16").Value = 8
